# Recomputed NATMI ligand-receptor scoring statistics using updated TPM expression values.
# Columns G:T (rows 2-17) are refreshed; columns A:F (identifiers / counts) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds the new values for columns G..T (in order) of one data row.
$newValues = @(
    @(19.28294533333333,57.848836,0.04564777115344932,0.04564777115344931,3,1,13.604331,40.812993,0.8107276168878804,0.8107276168878805,262.331570969572,2360.984138726148,0.0370079087234793,0.0370079087234793),
    @(19.28294533333333,57.848836,0.04564777115344932,0.04564777115344931,3,1,1.061748,3.185244,0.06327311690486458,0.06327311690486459,20.473628641776,184.262657775984,0.002888276760638704,0.002888276760638703),
    @(19.28294533333333,57.848836,0.04564777115344932,0.04564777115344931,3,1,1.995771333333333,5.987314,0.1189346934389115,0.1189346934389116,38.48434951850044,346.359145666504,0.005429103668305083,0.005429103668305083),
    @(19.28294533333333,57.848836,0.04564777115344932,0.04564777115344931,2,0.6666666666666666,0.1185463333333334,0.355639,0.007064572768343379,0.007064572768343379,2.285922465133778,20.573302186204,0.0003224820010262285,0.0003224820010262284),
    @(268.8003336666666,806.401001,0.6363206400827226,0.6363206400827226,3,1,13.604331,40.812993,0.8107276168878804,0.8107276168878805,3656.848712111776,32911.63840900599,0.5158827161108364,0.5158827161108365),
    @(268.8003336666666,806.401001,0.6363206400827226,0.6363206400827226,3,1,1.061748,3.185244,0.06327311690486458,0.06327311690486459,285.3982166699159,2568.583950029244,0.04026199024893237,0.04026199024893237),
    @(268.8003336666666,806.401001,0.6363206400827226,0.6363206400827226,3,1,1.995771333333333,5.987314,0.1189346934389115,0.1189346934389116,536.4640003223681,4828.176002901313,0.07568060025709057,0.07568060025709059),
    @(268.8003336666666,806.401001,0.6363206400827226,0.6363206400827226,2,0.6666666666666666,0.1185463333333334,0.355639,0.007064572768343379,0.007064572768343379,31.86529395495989,286.787645594639,0.004495333465863231,0.004495333465863231),
    @(56.43559133333333,169.306774,0.1335977939863952,0.1335977939863952,3,1,13.604331,40.812993,0.8107276168878804,0.8107276168878805,767.7684646793979,6909.916182114582,0.1083114211400682,0.1083114211400682),
    @(56.43559133333333,169.306774,0.1335977939863952,0.1335977939863952,3,1,1.061748,3.185244,0.06327311690486458,0.06327311690486459,59.92037622698399,539.2833860428559,0.008453148837133197,0.008453148837133198),
    @(56.43559133333333,169.306774,0.1335977939863952,0.1335977939863952,3,1,1.995771333333333,5.987314,0.1189346934389115,0.1189346934389116,112.6325353627818,1013.692818265036,0.01588941267188677,0.01588941267188677),
    @(56.43559133333333,169.306774,0.1335977939863952,0.1335977939863952,2,0.6666666666666666,0.1185463333333334,0.355639,0.007064572768343379,0.007064572768343379,6.690232422065112,60.212091798586,0.0009438113373070363,0.0009438113373070363),
    @(77.91019566666667,233.730587,0.184433794777433,0.1844337947774329,3,1,13.604331,40.812993,0.8107276168878804,0.8107276168878805,1059.916090124099,9539.244811116891,0.1495255709134966,0.1495255709134966),
    @(77.91019566666667,233.730587,0.184433794777433,0.1844337947774329,3,1,1.061748,3.185244,0.06327311690486458,0.06327311690486459,82.72099442869199,744.488949858228,0.01166970105816032,0.01166970105816032),
    @(77.91019566666667,233.730587,0.184433794777433,0.1844337947774329,3,1,1.995771333333333,5.987314,0.1189346934389115,0.1189346934389116,155.4909350859242,1399.418415773318,0.02193557684162911,0.02193557684162911),
    @(77.91019566666667,233.730587,0.184433794777433,0.1844337947774329,2,0.6666666666666666,0.1185463333333334,0.355639,0.007064572768343379,0.007064572768343379,9.23596802556589,83.12371223009302,0.001302945964146884,0.001302945964146884)
)

$startRow = 2
$startCol = 7  # column G

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowValues = $newValues[$i]
    $rowIndex = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($rowIndex, $startCol + $j).Value = $rowValues[$j]
    }
}

Write-Host "Updated TPM-derived statistics for rows 2-17 (columns G:T)."
